$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric stay as literal text (matches source formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.600.93'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '3.252.48'
$ws.Range("E3").Value = '  +2.98%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '608.52'
$ws.Range("E5").Value = '  +1.16%  '

$ws.Range("D6").Value = '157.96'
$ws.Range("E6").Value = '  +2.81%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.253.91'
$ws.Range("E8").Value = '  +3.15%  '

$ws.Range("D9").Value = '0.550'
$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  +2.88%  '

$ws.Range("D11").Value = '5.86'
$ws.Range("E11").Value = '  +6.59%  '

$ws.Range("D12").Value = '0.505'
$ws.Range("E12").Value = '  -0.83%  '

$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("D14").Value = '39.29'
$ws.Range("E14").Value = '  +2.67%  '

$ws.Range("D15").Value = '3.786.08'
$ws.Range("E15").Value = '  +2.92%  '

$ws.Range("D16").Value = '66.693.74'
$ws.Range("E16").Value = '  +0.73%  '

$ws.Range("D17").Value = '7.45'
$ws.Range("E17").Value = '  +1.21%  '

$ws.Range("D18").Value = '3.245.91'
$ws.Range("E18").Value = '  +2.71%  '

$ws.Range("E19").Value = '  +1.03%  '

$ws.Range("D20").Value = '508.21'
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").Value = '15.43'
$ws.Range("E21").Value = '  +0.49%  '

$ws.Range("D22").Value = '0.754'
$ws.Range("E22").Value = '  +3.87%  '

$ws.Range("D23").Value = '8.12'
$ws.Range("E23").Value = '  +0.60%  '

$ws.Range("E24").Value = '  +1.35%  '

$ws.Range("D25").Value = '86.90'
$ws.Range("E25").Value = '  +2.91%  '

$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '3.03'
$ws.Range("E27").Value = '  +1.63%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '9.12'
$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = '2.43'
$ws.Range("E29").Value = '  +1.92%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.135'
$ws.Range("E30").Value = '  +53.56%  '

$ws.Range("D31").Value = '2.91'
$ws.Range("E31").Value = '  -3.86%  '

$ws.Range("D32").Value = '6.86'
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("D33").Value = '28.12'
$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("E35").Value = '  -3.21%  '

$ws.Range("D36").Value = '6.47'
$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("E37").Value = '  +22.84%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '55.66'
$ws.Range("E38").Value = '  +1.70%  '

$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0787'
$ws.Range("E39").Value = '  +16.94%  '

$ws.Range("D40").Value = '494.77'
$ws.Range("E40").Value = '  -1.70%  '

$ws.Range("E41").Value = '  +2.07%  '

$ws.Range("E42").Value = '  +0.83%  '

$ws.Range("D43").Value = '8.86'
$ws.Range("E43").Value = '  +1.04%  '

$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.52'
$ws.Range("E44").Value = '  +4.92%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.293'
$ws.Range("E45").Value = '  -0.32%  '

$ws.Range("D46").Value = '2.981.70'
$ws.Range("E46").Value = '  +5.79%  '

$ws.Range("D47").Value = '29.11'
$ws.Range("E47").Value = '  +4.91%  '

$ws.Range("D48").Value = '2.51'
$ws.Range("E48").Value = '  +6.52%  '

$ws.Range("D49").Value = '0.120'
$ws.Range("E49").Value = '  +2.98%  '

$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").Value = '121.35'
$ws.Range("E51").Value = '  +0.31%  '
